# Apply updated numeric results to the workbook.
# The commit only changes computed values (re-run of the underlying
# statistical model); no structural changes are needed.

$wb = $excel.ActiveWorkbook

# --- Sheet "Full results" ---
$ws1 = $wb.Worksheets.Item("Full results")

$ws1.Range("H2").Value = 0.568327370344274
$ws1.Range("I2").Value = 0.182541134774655
$ws1.Range("O2").Value = 0.431724059767481

$ws1.Range("F3").Value = 0.578532758791308
$ws1.Range("G3").Value = 0.205295420524164

$ws1.Range("C4").Value = 0.637358737239101
$ws1.Range("D4").Value = 0.362731764845503
$ws1.Range("E4").Value = 1.0000905020846
$ws1.Range("J4").Value = 0.362698939830343
$ws1.Range("K4").Value = 0.205276842499685
$ws1.Range("L4").Value = 0.0102044649165169
$ws1.Range("M4").Value = 0.0690251199371374
$ws1.Range("N4").Value = 0.215481307416202

# --- Sheet "For plotting" ---
$ws2 = $wb.Worksheets.Item("For plotting")

$ws2.Range("C2").Value = 0.362698939830343
$ws2.Range("D2").Value = 0.328838568474349
$ws2.Range("E2").Value = 0.396559311186337

$ws2.Range("C3").Value = 0.215481307416202
$ws2.Range("D3").Value = 0.191412127455423
$ws2.Range("E3").Value = 0.239550487376981

$ws2.Range("C4").Value = 0.431724059767481
$ws2.Range("D4").Value = 0.399980320526992
$ws2.Range("E4").Value = 0.463467799007969
